$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value2 = 1.459612070389937
$ws.Range("C2").Value2 = 1.667794583268128
$ws.Range("D2").Value2 = 0.8054896365839992
$ws.Range("E2").Value2 = 0.496779210170732
$ws.Range("G2").Value2 = 4.429675500412797

$ws.Range("B3").Value2 = 1.459612070389937
$ws.Range("C3").Value2 = 1.667794583268128
$ws.Range("D3").Value2 = 0.1575252929769615
$ws.Range("E3").Value2 = 8.660232485948974
$ws.Range("G3").Value2 = 11.945164432584

$ws.Range("B4").Value2 = 1.459612070389937
$ws.Range("C4").Value2 = 1.667794583268128
$ws.Range("D4").Value2 = 0.8054896365839992
$ws.Range("E4").Value2 = 0.496779210170732
$ws.Range("G4").Value2 = 4.429675500412797

$ws.Range("B5").Value2 = 0.6753301551942219
$ws.Range("C5").Value2 = 1.667794583268128
$ws.Range("D5").Value2 = 0.8054896365839992
$ws.Range("E5").Value2 = 0.496779210170732
$ws.Range("G5").Value2 = 3.645393585217082

$ws.Range("B6").Value2 = 1.459612070389937
$ws.Range("C6").Value2 = 0.3127903958511391
$ws.Range("D6").Value2 = 0.1575252929769615
$ws.Range("E6").Value2 = 0.496779210170732
$ws.Range("G6").Value2 = 2.42670696938877

$ws.Range("B7").Value2 = 3.230985683306322
$ws.Range("C7").Value2 = 1.667794583268128
$ws.Range("D7").Value2 = 0.8054896365839992
$ws.Range("E7").Value2 = 0.496779210170732
$ws.Range("G7").Value2 = 6.201049113329182

$ws.Range("B8").Value2 = 0.6753301551942219
$ws.Range("C8").Value2 = 0.3127903958511391
$ws.Range("D8").Value2 = 0.1575252929769615
$ws.Range("E8").Value2 = 0.496779210170732
$ws.Range("G8").Value2 = 1.642425054193055

$ws.Range("B9").Value2 = 0.003994804209775715
$ws.Range("C9").Value2 = 0.002777888934908601
$ws.Range("D9").Value2 = 3.900430680208489
$ws.Range("E9").Value2 = 0.496779210170732
$ws.Range("G9").Value2 = 4.403982583523906

$ws.Range("B10").Value2 = 3.230985683306322
$ws.Range("C10").Value2 = 1.667794583268128
$ws.Range("D10").Value2 = 0.1575252929769615
$ws.Range("E10").Value2 = 0.496779210170732
$ws.Range("G10").Value2 = 5.553084769722144

$ws.Range("B11").Value2 = 3.230985683306322
$ws.Range("C11").Value2 = 1.667794583268128
$ws.Range("D11").Value2 = 0.8054896365839992
$ws.Range("E11").Value2 = 0.496779210170732
$ws.Range("G11").Value2 = 6.201049113329182

$ws.Range("B12").Value2 = 1.459612070389937
$ws.Range("C12").Value2 = 1.667794583268128
$ws.Range("D12").Value2 = 0.8054896365839992
$ws.Range("E12").Value2 = 0.496779210170732
$ws.Range("G12").Value2 = 4.429675500412797

$ws.Range("B13").Value2 = 3.230985683306322
$ws.Range("C13").Value2 = 1.667794583268128
$ws.Range("D13").Value2 = 0.8054896365839992
$ws.Range("E13").Value2 = 0.496779210170732
$ws.Range("G13").Value2 = 6.201049113329182

$ws.Range("B14").Value2 = 1.459612070389937
$ws.Range("C14").Value2 = 1.667794583268128
$ws.Range("D14").Value2 = 0.1575252929769615
$ws.Range("E14").Value2 = 0.496779210170732
$ws.Range("G14").Value2 = 3.781711156805759

$ws.Range("B15").Value2 = 0.127881588408715
$ws.Range("C15").Value2 = 0.3127903958511391
$ws.Range("D15").Value2 = 0.1575252929769615
$ws.Range("E15").Value2 = 0.496779210170732
$ws.Range("G15").Value2 = 1.094976487407548

$ws.Range("B16").Value2 = 3.230985683306322
$ws.Range("C16").Value2 = 1.667794583268128
$ws.Range("D16").Value2 = 0.1575252929769615
$ws.Range("E16").Value2 = 0.496779210170732
$ws.Range("G16").Value2 = 5.553084769722144

$ws.Range("B17").Value2 = 0.3048080303191223
$ws.Range("C17").Value2 = 0.00007097389502863649
$ws.Range("D17").Value2 = 0.8054896365839992
$ws.Range("E17").Value2 = 0.496779210170732
$ws.Range("G17").Value2 = 1.607147850968882

$ws.Range("B18").Value2 = 3.230985683306322
$ws.Range("C18").Value2 = 1.667794583268128
$ws.Range("D18").Value2 = 0.8054896365839992
$ws.Range("E18").Value2 = 0.496779210170732
$ws.Range("G18").Value2 = 6.201049113329182

$ws.Range("B19").Value2 = 0.6753301551942219
$ws.Range("C19").Value2 = 1.667794583268128
$ws.Range("D19").Value2 = 0.1575252929769615
$ws.Range("E19").Value2 = 8.660232485948974
$ws.Range("G19").Value2 = 11.16088251738829

$ws.Range("B20").Value2 = 3.230985683306322
$ws.Range("C20").Value2 = 1.667794583268128
$ws.Range("D20").Value2 = 26.21740644021617
$ws.Range("E20").Value2 = 0.496779210170732
$ws.Range("G20").Value2 = 31.61296591696135

$ws.Range("B21").Value2 = 3.230985683306322
$ws.Range("C21").Value2 = 1.667794583268128
$ws.Range("D21").Value2 = 3.900430680208489
$ws.Range("E21").Value2 = 8.660232485948974
$ws.Range("G21").Value2 = 17.45944343273191

$ws.Range("B22").Value2 = 3.230985683306322
$ws.Range("C22").Value2 = 1.667794583268128
$ws.Range("D22").Value2 = 0.1575252929769615
$ws.Range("E22").Value2 = 8.660232485948974
$ws.Range("G22").Value2 = 13.71653804550039

$ws.Range("B23").Value2 = 3.230985683306322
$ws.Range("C23").Value2 = 1.667794583268128
$ws.Range("D23").Value2 = 0.1575252929769615
$ws.Range("E23").Value2 = 0.496779210170732
$ws.Range("G23").Value2 = 5.553084769722144

$ws.Range("B24").Value2 = 1.459612070389937
$ws.Range("C24").Value2 = 1.667794583268128
$ws.Range("D24").Value2 = 0.8054896365839992
$ws.Range("E24").Value2 = 8.660232485948974
$ws.Range("G24").Value2 = 12.59312877619104

$ws.Range("B25").Value2 = 3.230985683306322
$ws.Range("C25").Value2 = 1.667794583268128
$ws.Range("D25").Value2 = 0.8054896365839992
$ws.Range("E25").Value2 = 0.496779210170732
$ws.Range("G25").Value2 = 6.201049113329182

$ws.Range("B26").Value2 = 0.127881588408715
$ws.Range("C26").Value2 = 0.3127903958511391
$ws.Range("D26").Value2 = 0.8054896365839992
$ws.Range("E26").Value2 = 0.496779210170732
$ws.Range("G26").Value2 = 1.742940831014585

$ws.Range("B27").Value2 = 0.3048080303191223
$ws.Range("C27").Value2 = 1.667794583268128
$ws.Range("D27").Value2 = 3.900430680208489
$ws.Range("E27").Value2 = 0.496779210170732
$ws.Range("G27").Value2 = 6.369812503966472

$ws.Range("B28").Value2 = 0.127881588408715
$ws.Range("C28").Value2 = 0.3127903958511391
$ws.Range("D28").Value2 = 3.900430680208489
$ws.Range("E28").Value2 = 8.660232485948974
$ws.Range("G28").Value2 = 13.00133515041732

$ws.Range("B29").Value2 = 3.230985683306322
$ws.Range("C29").Value2 = 1.667794583268128
$ws.Range("D29").Value2 = 0.1575252929769615
$ws.Range("E29").Value2 = 0.496779210170732
$ws.Range("G29").Value2 = 5.553084769722144

$ws.Range("B30").Value2 = 0.6753301551942219
$ws.Range("C30").Value2 = 1.667794583268128
$ws.Range("D30").Value2 = 3.900430680208489
$ws.Range("E30").Value2 = 0.496779210170732
$ws.Range("G30").Value2 = 6.740334628841572

$ws.Range("B31").Value2 = 0.6753301551942219
$ws.Range("C31").Value2 = 1.667794583268128
$ws.Range("D31").Value2 = 26.21740644021617
$ws.Range("E31").Value2 = 0.496779210170732
$ws.Range("G31").Value2 = 29.05731038884925

$ws.Range("B32").Value2 = 3.230985683306322
$ws.Range("C32").Value2 = 1.667794583268128
$ws.Range("D32").Value2 = 0.8054896365839992
$ws.Range("E32").Value2 = 0.496779210170732
$ws.Range("G32").Value2 = 6.201049113329182

$ws.Range("B33").Value2 = 0.127881588408715
$ws.Range("C33").Value2 = 1.667794583268128
$ws.Range("D33").Value2 = 0.1575252929769615
$ws.Range("E33").Value2 = 0.496779210170732
$ws.Range("G33").Value2 = 2.449980674824537

$ws.Range("B34").Value2 = 3.230985683306322
$ws.Range("C34").Value2 = 1.667794583268128
$ws.Range("D34").Value2 = 0.8054896365839992
$ws.Range("E34").Value2 = 8.660232485948974
$ws.Range("G34").Value2 = 14.36450238910742

$ws.Range("B35").Value2 = 0.01514828764759746
$ws.Range("C35").Value2 = 0.3127903958511391
$ws.Range("D35").Value2 = 3.900430680208489
$ws.Range("E35").Value2 = 0.496779210170732
$ws.Range("G35").Value2 = 4.725148573877958

$ws.Range("B36").Value2 = 0.6753301551942219
$ws.Range("C36").Value2 = 1.667794583268128
$ws.Range("D36").Value2 = 0.1575252929769615
$ws.Range("E36").Value2 = 0.496779210170732
$ws.Range("G36").Value2 = 2.997429241610044

$ws.Range("B37").Value2 = 1.459612070389937
$ws.Range("C37").Value2 = 114.8270160096505
$ws.Range("D37").Value2 = 0.8054896365839992
$ws.Range("E37").Value2 = 8.660232485948974
$ws.Range("G37").Value2 = 125.7523502025734

$ws.Range("B38").Value2 = 0.127881588408715
$ws.Range("C38").Value2 = 0.3127903958511391
$ws.Range("D38").Value2 = 0.1575252929769615
$ws.Range("E38").Value2 = 0.496779210170732
$ws.Range("G38").Value2 = 1.094976487407548
